# Adds two more rows of "test" login data to Sheet1, marks the header
# columns on Sheet2 (Adults/Children per Room) as wrap-text, and merges
# Sheet2's hotel-booking data into Sheet3 alongside the existing billing
# data (shifted right), mirroring the "Added Log4j and ExtentReport" test
# data consolidation commit.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------------
# Sheet1: clear row 3, change B4, append two new "test" rows (5 and 6)
# ---------------------------------------------------------------------
$ws1.Range("A3").Value = ""
$ws1.Range("B3").Value = ""

$ws1.Range("B4").Value = "test"

# Row 2 carries the style (s="1") we need on the new rows - copy its
# format down before overwriting the values.
$ws1.Range("A2:B2").Copy()
$ws1.Range("A5:B6").PasteSpecial(-4122)

$ws1.Range("A5").Value = "test"
$ws1.Range("B5").Value = "JayanthiSaraPeru@123"
$ws1.Range("A6").Value = "test"
$ws1.Range("B6").Value = "test"

$ws1.Columns.Item(1).ColumnWidth = 19.072916666666668
$ws1.Columns.Item(2).ColumnWidth = 20.893229166666668

$ws1.Range("B6").Select()

# ---------------------------------------------------------------------
# Sheet2: wrap the last two header cells, resize columns
# ---------------------------------------------------------------------
$ws2.Range("G1:H1").WrapText = $true

$ws2.Columns.Item(1).ColumnWidth = 10.619791666666666
$ws2.Columns.Item(2).ColumnWidth = 13.529947916666666
$ws2.Columns.Item(3).ColumnWidth = 15.983072916666666
$ws2.Columns.Item(4).ColumnWidth = 15.799479166666666
$ws2.Columns.Item(5).ColumnWidth = 15.256510416666666
$ws2.Columns.Item(6).ColumnWidth = 16.619791666666668
$ws2.Columns.Item(7).ColumnWidth = 15.436197916666666
$ws2.Columns.Item(8).ColumnWidth = 18.436197916666668

$ws2.Range("A1:H2").Select()

# ---------------------------------------------------------------------
# Sheet3: shift the existing billing-info table 8 columns to the right
# (I:P) and paste Sheet2's booking table into the freed A:H columns.
# ---------------------------------------------------------------------
$ws3.Range("A1:H2").Insert(-4161)

$ws2.Range("A1:H2").Copy()
$ws3.Range("A1").PasteSpecial(-4163)
$ws2.Range("A1:H2").Copy()
$ws3.Range("A1").PasteSpecial(-4122)

$ws3.Columns.Item(4).ColumnWidth = 13.256510416666666
$ws3.Columns.Item(5).ColumnWidth = 12.799479166666666
$ws3.Columns.Item(6).ColumnWidth = 16.893229166666668
$ws3.Columns.Item(7).ColumnWidth = 11.709635416666666
$ws3.Columns.Item(12).ColumnWidth = 19.072916666666668

$ws3.Range("F6").Select()
